# Fix: prevent hidden columns from being labeled upon detecting changes (#11)
#
# The worksheet contains a FV2504-vs-FV2410 comparison table. Column L
# ("Änderung") was previously flagging a number of rows as "ÄNDERUNG"
# (changed) even though the only differences were in columns that should
# have been ignored (hidden columns). After the fix, those rows no longer
# carry the "ÄNDERUNG" label: the L cell is cleared and reformatted.
#
# In addition, whenever such a row is also the first row of a new
# "Segmentname" group (column B), the whole row is re-styled to the
# "group header" look (the same look already used by the pre-existing
# group-header rows, e.g. row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cells that already carry the target styles (row 2 is an
# existing, untouched "group header" row):
#   A2, C2:K2, M2:V2 -> style used for most cells of a header row
#   B2               -> style used for column B of a header row
#   L2               -> style used for a cleared "Änderung" cell
$tmplA   = $ws.Range("A2")
$tmplB   = $ws.Range("B2")
$tmplCK  = $ws.Range("C2:K2")
$tmplMV  = $ws.Range("M2:V2")
$tmplL   = $ws.Range("L2")

# Rows whose "Änderung" (column L) flag must be cleared because the
# previously-detected difference only existed in a column that should
# have been skipped (hidden column false positive).
$clearRows = @(
    38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,54,55,56,57,58,59,60,
    61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79,80,81,82,
    83,84,85,86,87,88,89,90,91,92,93,94,96,97,98,99,100,101,102,103,
    104,105,107,108,109,110,111,112,113,114,115,116,117,118,119,120,
    121,122,123,124,125,126,127,128,129,130,131,132,133,134,135,136,
    137,138,139,140,141,142,143,144,145,146,147,148,149,150,151,152,
    153,154,155,156,157,158,159,160,161,162,163,164,165,166,167,168,
    169,170,171,172,173
)

# Subset of the rows above that are also the first row of a new
# "Segmentname" group, and therefore get the full group-header styling.
$headerRows = @(38,42,48,52,59,65,67,70,74,78,82,85,89,100,104,108,117,123,145,168,171)
$headerRowSet = @{}
foreach ($hr in $headerRows) { $headerRowSet[$hr] = $true }

foreach ($r in $clearRows) {
    if ($headerRowSet.ContainsKey($r)) {
        $tmplA.Copy()
        $ws.Range("A$r").PasteSpecial(-4122)

        $tmplB.Copy()
        $ws.Range("B$r").PasteSpecial(-4122)

        $tmplCK.Copy()
        $ws.Range("C$r`:K$r").PasteSpecial(-4122)

        $tmplMV.Copy()
        $ws.Range("M$r`:V$r").PasteSpecial(-4122)
    }

    # Clear the "Änderung" label and restore the (now-unflagged) style.
    $lcell = $ws.Range("L$r")
    $lcell.ClearContents()
    $tmplL.Copy()
    $lcell.PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
